# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right after "总计", populated with the
#    Q3 fund holdings table (built by duplicating the "2022-Q2" sheet so the
#    sheetPr / column formatting / header styling come along for free, then
#    overwriting the data rows).
# 2) Insert a new second row into "总计" summarizing the 2022-Q3 quarter,
#    pushing the existing quarters (2022-Q2, 2022-Q1, 2021-Q4) down by one
#    row.
#
# Text-looking values (fund codes with leading zeros, decimal numbers stored
# as text in the source data, …) are written through a formula + paste-values
# round trip so they land as real text cells instead of being auto-coerced to
# numbers by plain `.Value =` assignment (which would e.g. turn "008185"
# into 8185).

function Set-TextValue($cell, [string]$text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- 1) Build the new "2022-Q3" sheet -------------------------------------
$q2Sheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# The duplicated sheet has 6 data rows (rows 2-7); the Q3 table only needs 4
# (rows 2-5), so drop the trailing two rows.
$q3Sheet.Rows.Item(7).Delete()
$q3Sheet.Rows.Item(6).Delete()

# code, name, size, stock position, position pct, market value, rank
$q3Data = @(
    @("008185", "诺安研究优选混合A", "1.36", "93.89", "6.60", "0.0898", 3),
    @("007316", "交银施罗德可转债债券A", "0.66", "23.69", "0.84", "0.0055", 7),
    @("014497", "诺安研究优选混合C", "0.08", "93.89", "6.60", "0.0053", 3),
    @("007317", "交银施罗德可转债债券C", "0.38", "23.69", "0.84", "0.0032", 7)
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $r = $i + 2
    $row = $q3Data[$i]
    $q3Sheet.Cells.Item($r, 1).Value = $i
    Set-TextValue $q3Sheet.Cells.Item($r, 2) $row[0]
    Set-TextValue $q3Sheet.Cells.Item($r, 3) $row[1]
    Set-TextValue $q3Sheet.Cells.Item($r, 4) $row[2]
    Set-TextValue $q3Sheet.Cells.Item($r, 5) $row[3]
    Set-TextValue $q3Sheet.Cells.Item($r, 6) $row[4]
    Set-TextValue $q3Sheet.Cells.Item($r, 7) $row[5]
    $q3Sheet.Cells.Item($r, 8).Value = $row[6]
}

# --- 2) Insert the 2022-Q3 row into "总计" ---------------------------------
# Shift the existing quarter rows (2022-Q2, 2022-Q1, 2021-Q4) down one row by
# copying bottom-up (keeps formatting / styles, no new style entries), then
# write the new 2022-Q3 summary row into the now-empty row 2 and fix up the
# row-index column (A).
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

$totalSheet.Cells.Item(2, 1).Value = 0
Set-TextValue $totalSheet.Cells.Item(2, 2) "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 0.1

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
